$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 0.8029023333333333
$ws.Range("H2").Value = 2.408707
$ws.Range("I2").Value = 0.06206726394886004
$ws.Range("J2").Value = 0.06206726394886004
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 34.218763
$ws.Range("N2").Value = 102.656289
$ws.Range("O2").Value = 0.4046921425624349
$ws.Range("P2").Value = 0.4046921425624349
$ws.Range("Q2").Value = 27.47432465648033
$ws.Range("R2").Value = 247.268921908323
$ws.Range("S2").Value = 0.02511813403045235
$ws.Range("T2").Value = 0.02511813403045234

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 0.8029023333333333
$ws.Range("H3").Value = 2.408707
$ws.Range("I3").Value = 0.06206726394886004
$ws.Range("J3").Value = 0.06206726394886004
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 13.95683033333333
$ws.Range("N3").Value = 41.870491
$ws.Range("O3").Value = 0.165062061740135
$ws.Range("P3").Value = 0.165062061740135
$ws.Range("Q3").Value = 11.20597164057078
$ws.Range("R3").Value = 100.853744765137
$ws.Range("S3").Value = 0.01024495055396799
$ws.Range("T3").Value = 0.01024495055396799

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 0.8029023333333333
$ws.Range("H4").Value = 2.408707
$ws.Range("I4").Value = 0.06206726394886004
$ws.Range("J4").Value = 0.06206726394886004
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 14.45863566666667
$ws.Range("N4").Value = 43.375907
$ws.Range("O4").Value = 0.1709967203219173
$ws.Range("P4").Value = 0.1709967203219173
$ws.Range("Q4").Value = 11.60887231358322
$ws.Range("R4").Value = 104.479850822249
$ws.Range("S4").Value = 0.01061329857460984
$ws.Range("T4").Value = 0.01061329857460984

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 0.8029023333333333
$ws.Range("H5").Value = 2.408707
$ws.Range("I5").Value = 0.06206726394886004
$ws.Range("J5").Value = 0.06206726394886004
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 10.640006
$ws.Range("N5").Value = 31.920018
$ws.Range("O5").Value = 0.1258352566694817
$ws.Range("P5").Value = 0.1258352566694817
$ws.Range("Q5").Value = 8.542885644080668
$ws.Range("R5").Value = 76.885970796726
$ws.Range("S5").Value = 0.00781025008977727
$ws.Range("T5").Value = 0.00781025008977727

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 0.8029023333333333
$ws.Range("H6").Value = 2.408707
$ws.Range("I6").Value = 0.06206726394886004
$ws.Range("J6").Value = 0.06206726394886004
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 8.185362
$ws.Range("N6").Value = 24.556086
$ws.Range("O6").Value = 0.09680512663269379
$ws.Range("P6").Value = 0.09680512663269379
$ws.Range("Q6").Value = 6.572046248977999
$ws.Range("R6").Value = 59.14841624080199
$ws.Range("S6").Value = 0.006008429346314226
$ws.Range("T6").Value = 0.006008429346314226

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 0.8029023333333333
$ws.Range("H7").Value = 2.408707
$ws.Range("I7").Value = 0.06206726394886004
$ws.Range("J7").Value = 0.06206726394886004
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 3.095449666666667
$ws.Range("N7").Value = 9.286349
$ws.Range("O7").Value = 0.03660869207333731
$ws.Range("P7").Value = 0.03660869207333731
$ws.Range("Q7").Value = 2.485343760082555
$ws.Range("R7").Value = 22.368093840743
$ws.Range("S7").Value = 0.002272201353738367
$ws.Range("T7").Value = 0.002272201353738367

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 2.159929333333333
$ws.Range("H8").Value = 6.479788
$ws.Range("I8").Value = 0.1669703754456877
$ws.Range("J8").Value = 0.1669703754456877
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 34.218763
$ws.Range("N8").Value = 102.656289
$ws.Range("O8").Value = 0.4046921425624349
$ws.Range("P8").Value = 0.4046921425624349
$ws.Range("Q8").Value = 73.91010995408134
$ws.Range("R8").Value = 665.190989586732
$ws.Range("S8").Value = 0.06757159898356951
$ws.Range("T8").Value = 0.0675715989835695

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 2.159929333333333
$ws.Range("H9").Value = 6.479788
$ws.Range("I9").Value = 0.1669703754456877
$ws.Range("J9").Value = 0.1669703754456877
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 13.95683033333333
$ws.Range("N9").Value = 41.870491
$ws.Range("O9").Value = 0.165062061740135
$ws.Range("P9").Value = 0.165062061740135
$ws.Range("Q9").Value = 30.14576723732311
$ws.Range("R9").Value = 271.311905135908
$ws.Range("S9").Value = 0.02756047442058962
$ws.Range("T9").Value = 0.02756047442058962

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 2.159929333333333
$ws.Range("H10").Value = 6.479788
$ws.Range("I10").Value = 0.1669703754456877
$ws.Range("J10").Value = 0.1669703754456877
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 14.45863566666667
$ws.Range("N10").Value = 43.375907
$ws.Range("O10").Value = 0.1709967203219173
$ws.Range("P10").Value = 0.1709967203219173
$ws.Range("Q10").Value = 31.22963129641289
$ws.Range("R10").Value = 281.066681667716
$ws.Range("S10").Value = 0.02855138659213178
$ws.Range("T10").Value = 0.02855138659213178

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 2.159929333333333
$ws.Range("H11").Value = 6.479788
$ws.Range("I11").Value = 0.1669703754456877
$ws.Range("J11").Value = 0.1669703754456877
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 10.640006
$ws.Range("N11").Value = 31.920018
$ws.Range("O11").Value = 0.1258352566694817
$ws.Range("P11").Value = 0.1258352566694817
$ws.Range("Q11").Value = 22.98166106624267
$ws.Range("R11").Value = 206.834949596184
$ws.Range("S11").Value = 0.02101076005040783
$ws.Range("T11").Value = 0.02101076005040783

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 2.159929333333333
$ws.Range("H12").Value = 6.479788
$ws.Range("I12").Value = 0.1669703754456877
$ws.Range("J12").Value = 0.1669703754456877
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 8.185362
$ws.Range("N12").Value = 24.556086
$ws.Range("O12").Value = 0.09680512663269379
$ws.Range("P12").Value = 0.09680512663269379
$ws.Range("Q12").Value = 17.679803487752
$ws.Range("R12").Value = 159.118231389768
$ws.Range("S12").Value = 0.01616358833892822
$ws.Range("T12").Value = 0.01616358833892822

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 2.159929333333333
$ws.Range("H13").Value = 6.479788
$ws.Range("I13").Value = 0.1669703754456877
$ws.Range("J13").Value = 0.1669703754456877
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 3.095449666666667
$ws.Range("N13").Value = 9.286349
$ws.Range("O13").Value = 0.03660869207333731
$ws.Range("P13").Value = 0.03660869207333731
$ws.Range("Q13").Value = 6.685952534890222
$ws.Range("R13").Value = 60.173572814012
$ws.Range("S13").Value = 0.006112567060060701
$ws.Range("T13").Value = 0.006112567060060701

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 9.973171666666666
$ws.Range("H14").Value = 29.919515
$ws.Range("I14").Value = 0.7709623606054523
$ws.Range("J14").Value = 0.7709623606054523
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 34.218763
$ws.Range("N14").Value = 102.656289
$ws.Range("O14").Value = 0.4046921425624349
$ws.Range("P14").Value = 0.4046921425624349
$ws.Range("Q14").Value = 341.2695976199817
$ws.Range("R14").Value = 3071.426378579835
$ws.Range("S14").Value = 0.3120024095484131
$ws.Range("T14").Value = 0.312002409548413

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 9.973171666666666
$ws.Range("H15").Value = 29.919515
$ws.Range("I15").Value = 0.7709623606054523
$ws.Range("J15").Value = 0.7709623606054523
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 13.95683033333333
$ws.Range("N15").Value = 41.870491
$ws.Range("O15").Value = 0.165062061740135
$ws.Range("P15").Value = 0.165062061740135
$ws.Range("Q15").Value = 139.1938648368739
$ws.Range("R15").Value = 1252.744783531865
$ws.Range("S15").Value = 0.1272566367655774
$ws.Range("T15").Value = 0.1272566367655774

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 9.973171666666666
$ws.Range("H16").Value = 29.919515
$ws.Range("I16").Value = 0.7709623606054523
$ws.Range("J16").Value = 0.7709623606054523
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 14.45863566666667
$ws.Range("N16").Value = 43.375907
$ws.Range("O16").Value = 0.1709967203219173
$ws.Range("P16").Value = 0.1709967203219173
$ws.Range("Q16").Value = 144.1984555694561
$ws.Range("R16").Value = 1297.786100125105
$ws.Range("S16").Value = 0.1318320351551757
$ws.Range("T16").Value = 0.1318320351551757

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 9.973171666666666
$ws.Range("H17").Value = 29.919515
$ws.Range("I17").Value = 0.7709623606054523
$ws.Range("J17").Value = 0.7709623606054523
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 10.640006
$ws.Range("N17").Value = 31.920018
$ws.Range("O17").Value = 0.1258352566694817
$ws.Range("P17").Value = 0.1258352566694817
$ws.Range("Q17").Value = 106.1146063723633
$ws.Range("R17").Value = 955.03145735127
$ws.Range("S17").Value = 0.09701424652929658
$ws.Range("T17").Value = 0.09701424652929658

$ws.Range("E18").Value = 3
$ws.Range("G18").Value = 9.973171666666666
$ws.Range("H18").Value = 29.919515
$ws.Range("I18").Value = 0.7709623606054523
$ws.Range("J18").Value = 0.7709623606054523
$ws.Range("K18").Value = 3
$ws.Range("M18").Value = 8.185362
$ws.Range("N18").Value = 24.556086
$ws.Range("O18").Value = 0.09680512663269379
$ws.Range("P18").Value = 0.09680512663269379
$ws.Range("Q18").Value = 81.63402037980998
$ws.Range("R18").Value = 734.7061834182899
$ws.Range("S18").Value = 0.07463310894745134
$ws.Range("T18").Value = 0.07463310894745134

$ws.Range("E19").Value = 3
$ws.Range("G19").Value = 9.973171666666666
$ws.Range("H19").Value = 29.919515
$ws.Range("I19").Value = 0.7709623606054523
$ws.Range("J19").Value = 0.7709623606054523
$ws.Range("K19").Value = 3
$ws.Range("M19").Value = 3.095449666666667
$ws.Range("N19").Value = 9.286349
$ws.Range("O19").Value = 0.03660869207333731
$ws.Range("P19").Value = 0.03660869207333731
$ws.Range("Q19").Value = 30.87145091119277
$ws.Range("R19").Value = 277.843058200735
$ws.Range("S19").Value = 0.02822392365953824
$ws.Range("T19").Value = 0.02822392365953824

